# Fruta / hortaliza, semanal
# Insert 3 new daily-price rows for Vega Monumental Concepción - Limón
# right before the existing row 384, shifting all subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at position 384 (old rows 384.. shift down to 387..)
$ws.Rows("384:386").Insert()

# Common (constant across the whole sheet) column values
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$prodId    = 100102
$producto  = "Cítricos"
$catId     = 100102003
$categoria = "Limón"
$variedad  = "Sin especificar"
$unidad    = "$/malla 16 kilos"
$kgUnidad  = 16

# Row 384
$ws.Range("A384").Value = $mercadoId
$ws.Range("B384").Value = $mercado
$ws.Range("C384").Value = $region
$ws.Range("D384").Value = 44694
$ws.Range("E384").Value = $codreg
$ws.Range("F384").Value = $tipo
$ws.Range("G384").Value = $prodId
$ws.Range("H384").Value = $producto
$ws.Range("I384").Value = $catId
$ws.Range("J384").Value = $categoria
$ws.Range("K384").Value = $variedad
$ws.Range("L384").Value = "1a amarillo"
$ws.Range("M384").Value = 450
$ws.Range("N384").Value = 14000
$ws.Range("O384").Value = 15000
$ws.Range("P384").Value = 14556
$ws.Range("Q384").Value = $unidad
$ws.Range("R384").Value = "Provincia de Melipilla"
$ws.Range("S384").Value = 910
$ws.Range("T384").Value = $kgUnidad

# Row 385
$ws.Range("A385").Value = $mercadoId
$ws.Range("B385").Value = $mercado
$ws.Range("C385").Value = $region
$ws.Range("D385").Value = 44694
$ws.Range("E385").Value = $codreg
$ws.Range("F385").Value = $tipo
$ws.Range("G385").Value = $prodId
$ws.Range("H385").Value = $producto
$ws.Range("I385").Value = $catId
$ws.Range("J385").Value = $categoria
$ws.Range("K385").Value = $variedad
$ws.Range("L385").Value = "1a plateado"
$ws.Range("M385").Value = 350
$ws.Range("N385").Value = 14000
$ws.Range("O385").Value = 15000
$ws.Range("P385").Value = 14429
$ws.Range("Q385").Value = $unidad
$ws.Range("R385").Value = "Provincia de Melipilla"
$ws.Range("S385").Value = 902
$ws.Range("T385").Value = $kgUnidad

# Row 386
$ws.Range("A386").Value = $mercadoId
$ws.Range("B386").Value = $mercado
$ws.Range("C386").Value = $region
$ws.Range("D386").Value = 44694
$ws.Range("E386").Value = $codreg
$ws.Range("F386").Value = $tipo
$ws.Range("G386").Value = $prodId
$ws.Range("H386").Value = $producto
$ws.Range("I386").Value = $catId
$ws.Range("J386").Value = $categoria
$ws.Range("K386").Value = $variedad
$ws.Range("L386").Value = "2a amarillo"
$ws.Range("M386").Value = 210
$ws.Range("N386").Value = 11000
$ws.Range("O386").Value = 12000
$ws.Range("P386").Value = 11524
$ws.Range("Q386").Value = $unidad
$ws.Range("R386").Value = "Provincia de Melipilla"
$ws.Range("S386").Value = 720
$ws.Range("T386").Value = $kgUnidad
